$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the remaining data rows (3-5) to their new values.
$ws.Range("A3").Value = "a and b"
$ws.Range("M3").Value = 0

$ws.Range("A4").Value = "c"

$ws.Range("A5").Value = "פרסון ואורח\ת"
$ws.Range("M5").Value = 1

# Remove the now-unused trailing rows (6-9).
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()
$ws.Rows.Item(6).Delete()
